{"js": "// Add three new paragraphs (\"Pove pora\", \"ETV Jabardasth\", \"Dhee Jodi\")\n// at the end of the document body, after the existing last paragraph (\"O\").\nconst body = context.document.body;\n\nconst newLines = [\"Pove pora\", \"ETV Jabardasth\", \"Dhee Jodi\"];\nfor (const line of newLines) {\n  body.insertParagraph(line, Word.InsertLocation.end);\n}\n\nawait context.sync();\n", "ps1": "# Add three new paragraphs (\"Pove pora\", \"ETV Jabardasth\", \"Dhee Jodi\")\n# at the end of the document, after the existing last paragraph (\"O\").\n$d = $word.ActiveDocument\n\n$lines = @(\"Pove pora\", \"ETV Jabardasth\", \"Dhee Jodi\")\nforeach ($line in $lines) {\n    $lastPara = $d.Paragraphs.Last\n    $lastPara.Range.InsertParagraphAfter()\n    $d.Paragraphs.Last.Range.Text = $line\n}\n"}
